$d = $word.ActiveDocument

$pairs = @(
    @("94×31=2914", "41×60=2460"),
    @("39×24=936", "22×51=1122"),
    @("57×41=2337", "89×72=6408"),
    @("27×82=2214", "69×78=5382"),
    @("80×20=1600", "15×29=435"),
    @("23×97=2231", "96×33=3168"),
    @("95×89=8455", "41×79=3239"),
    @("38×77=2926", "19×66=1254"),
    @("18×84=1512", "53×76=4028"),
    @("79×91=7189", "15×16=240"),
    @("68×28=1904", "85×22=1870"),
    @("65×38=2470", "37×47=1739"),
    @("34×95=3230", "65×84=5460"),
    @("68×93=6324", "25×26=650"),
    @("63×60=3780", "36×98=3528"),
    @("63×17=1071", "63×86=5418"),
    @("44×41=1804", "21×21=441"),
    @("95×52=4940", "77×39=3003"),
    @("53×69=3657", "82×28=2296"),
    @("34×86=2924", "84×90=7560"),
    @("53×73=3869", "58×17=986"),
    @("91×96=8736", "49×71=3479"),
    @("55×24=1320", "49×96=4704"),
    @("49×60=2940", "33×55=1815"),
    @("88×73=6424", "50×72=3600")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
